$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.036335610413886
$ws.Range("D2").Value = 1.036707749220278
$ws.Range("E2").Value = 1.053338013818284
$ws.Range("F2").Value = 1.059300192726017
$ws.Range("I2").Value = 1.033245217137009
$ws.Range("J2").Value = 1.041444486482763
$ws.Range("K2").Value = 1.0395007198573
$ws.Range("L2").Value = 1.056084258852767
$ws.Range("M2").Value = 1.062030073017463
$ws.Range("N2").Value = 1.042923457108921

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037948973991986
$ws.Range("D3").Value = 1.037883026318039
$ws.Range("E3").Value = 1.054940906084475
$ws.Range("F3").Value = 1.061061365101226
$ws.Range("I3").Value = 1.033621730457677
$ws.Range("J3").Value = 1.042698429936646
$ws.Range("K3").Value = 1.040484684890655
$ws.Range("L3").Value = 1.057498153647113
$ws.Range("M3").Value = 1.063603053022115
$ws.Range("N3").Value = 1.044179181306337

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038990042956107
$ws.Range("D4").Value = 1.038640820074643
$ws.Range("E4").Value = 1.055976589284304
$ws.Range("F4").Value = 1.062199395332664
$ws.Range("I4").Value = 1.033862425593763
$ws.Range("J4").Value = 1.043506623500186
$ws.Range("K4").Value = 1.041118124702668
$ws.Range("L4").Value = 1.058411063236324
$ws.Range("M4").Value = 1.064618861336648
$ws.Range("N4").Value = 1.044988522597439

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039427029632583
$ws.Range("D5").Value = 1.038958760392957
$ws.Range("E5").Value = 1.056411643722245
$ws.Range("F5").Value = 1.06267746067458
$ws.Range("I5").Value = 1.033962914604342
$ws.Range("J5").Value = 1.043845633179999
$ws.Range("K5").Value = 1.041383651545404
$ws.Range("L5").Value = 1.058794387971641
$ws.Range("M5").Value = 1.065045436079183
$ws.Range("N5").Value = 1.045328013709883

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039500362164619
$ws.Range("D6").Value = 1.03901210688501
$ws.Range("E6").Value = 1.056484671265158
$ws.Range("F6").Value = 1.062757709045087
$ws.Range("I6").Value = 1.033979746241228
$ws.Range("J6").Value = 1.043902510395298
$ws.Range("K6").Value = 1.041428189649567
$ws.Range("L6").Value = 1.058858723061695
$ws.Range("M6").Value = 1.065117032558397
$ws.Range("N6").Value = 1.04538497169735

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038995884643816
$ws.Range("D7").Value = 1.038645070897897
$ws.Range("E7").Value = 1.05598240385017
$ws.Range("F7").Value = 1.062205784672785
$ws.Range("I7").Value = 1.033863771075804
$ws.Range("J7").Value = 1.043511156317597
$ws.Range("K7").Value = 1.041121675705933
$ws.Range("L7").Value = 1.058416187043083
$ws.Range("M7").Value = 1.064624563084642
$ws.Range("N7").Value = 1.04499306185197

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036881458831213
$ws.Range("D8").Value = 1.037105501017158
$ws.Range("E8").Value = 1.053880035132238
$ws.Range("F8").Value = 1.059895720900379
$ws.Range("I8").Value = 1.033373071063797
$ws.Range("J8").Value = 1.041868928974059
$ws.Range("K8").Value = 1.039833932976304
$ws.Range("L8").Value = 1.056562506374704
$ws.Range("M8").Value = 1.062562093009496
$ws.Range("N8").Value = 1.043348502357238

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.033132932485288
$ws.Range("D9").Value = 1.034371643023869
$ws.Range("E9").Value = 1.050163437123736
$ws.Range("F9").Value = 1.055812543241163
$ws.Range("I9").Value = 1.032485776750365
$ws.Range("J9").Value = 1.038950254263901
$ws.Range("K9").Value = 1.037539541877318
$ws.Range("L9").Value = 1.053280493782769
$ws.Range("M9").Value = 1.058911816211804
$ws.Range("N9").Value = 1.040425682794216

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030617933222223
$ws.Range("D10").Value = 1.032534507281851
$ws.Range("E10").Value = 1.047676960072453
$ws.Range("F10").Value = 1.053081206383013
$ws.Range("I10").Value = 1.031878832833938
$ws.Range("J10").Value = 1.036987170594801
$ws.Range("K10").Value = 1.035992533519511
$ws.Range("L10").Value = 1.051081350540496
$ws.Range("M10").Value = 1.056466865732703
$ws.Range("N10").Value = 1.038459811321137

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029524949407738
$ws.Range("D11").Value = 1.031735443485938
$ws.Range("E11").Value = 1.046598055911894
$ws.Range("F11").Value = 1.051896146386361
$ws.Range("I11").Value = 1.031612316348034
$ws.Range("J11").Value = 1.036132898542676
$ws.Range("K11").Value = 1.035318428794749
$ws.Range("L11").Value = 1.050126310291619
$ws.Range("M11").Value = 1.055405303620334
$ws.Range("N11").Value = 1.037604326104711

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029118354892474
$ws.Range("D12").Value = 1.031438089413927
$ws.Range("E12").Value = 1.046196952182751
$ws.Range("F12").Value = 1.051455590696645
$ws.Range("I12").Value = 1.031512759359097
$ws.Range("J12").Value = 1.03581493479352
$ws.Range("K12").Value = 1.035067390315669
$ws.Range("L12").Value = 1.04977113300599
$ws.Range("M12").Value = 1.055010545276829
$ws.Range("N12").Value = 1.037285910810561

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029205598715121
$ws.Range("D13").Value = 1.031501897743273
$ws.Range("E13").Value = 1.046283006384439
$ws.Range("F13").Value = 1.051550108449394
$ws.Range("I13").Value = 1.031534140120306
$ws.Range("J13").Value = 1.035883168654718
$ws.Range("K13").Value = 1.035121268298359
$ws.Range("L13").Value = 1.049847339529438
$ws.Range("M13").Value = 1.05509524272912
$ws.Range("N13").Value = 1.037354241571668

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029491352742392
$ws.Range("D14").Value = 1.031710875303655
$ws.Range("E14").Value = 1.046564907785972
$ws.Range("F14").Value = 1.051859737584957
$ws.Range("I14").Value = 1.031604098402671
$ws.Range("J14").Value = 1.036106628867766
$ws.Range("K14").Value = 1.035297691105663
$ws.Range("L14").Value = 1.050096960129229
$ws.Range("M14").Value = 1.055372681969996
$ws.Range("N14").Value = 1.03757801912385

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029667333731932
$ws.Range("D15").Value = 1.031839560610782
$ws.Range("E15").Value = 1.046738549662802
$ws.Range("F15").Value = 1.052050460613601
$ws.Range("I15").Value = 1.031647127559561
$ws.Range("J15").Value = 1.03624422372595
$ws.Range("K15").Value = 1.035406305152262
$ws.Range("L15").Value = 1.050250701886959
$ws.Range("M15").Value = 1.055543561806169
$ws.Range("N15").Value = 1.037715809382515

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030690385999362
$ws.Range("D16").Value = 1.0325874624451
$ws.Range("E16").Value = 1.047748514962983
$ws.Range("F16").Value = 1.053159803627218
$ws.Range("I16").Value = 1.031896442226254
$ws.Range("J16").Value = 1.037043775427306
$ws.Range("K16").Value = 1.036037181536906
$ws.Range("L16").Value = 1.051144673416925
$ws.Range("M16").Value = 1.056537256205042
$ws.Range("N16").Value = 1.038516496538996

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.031331045678612
$ws.Range("D17").Value = 1.033055638495014
$ws.Range("E17").Value = 1.048381429337389
$ws.Range("F17").Value = 1.053855019223091
$ws.Range("I17").Value = 1.032051835679836
$ws.Range("J17").Value = 1.037544168820406
$ws.Range("K17").Value = 1.03643177218528
$ws.Range("L17").Value = 1.051704680495484
$ws.Range("M17").Value = 1.05715979300366
$ws.Range("N17").Value = 1.039017600548111

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.031704349349039
$ws.Range("D18").Value = 1.033328373680561
$ws.Range("E18").Value = 1.048750382079084
$ws.Range("F18").Value = 1.054260298775509
$ws.Range("I18").Value = 1.032142116783445
$ws.Range("J18").Value = 1.037835631131653
$ws.Range("K18").Value = 1.036661521758398
$ws.Range("L18").Value = 1.05203105416775
$ws.Range("M18").Value = 1.057522631017408
$ws.Range("N18").Value = 1.03930947676927

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031831571858127
$ws.Range("D19").Value = 1.033421311272417
$ws.Range("E19").Value = 1.04887614928759
$ws.Range("F19").Value = 1.054398450421223
$ws.Range("I19").Value = 1.032172839832236
$ws.Range("J19").Value = 1.037934943272024
$ws.Range("K19").Value = 1.036739791410487
$ws.Range("L19").Value = 1.052142293912172
$ws.Range("M19").Value = 1.05764630280798
$ws.Range("N19").Value = 1.039408929944272

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.031262348606441
$ws.Range("D20").Value = 1.033005443284532
$ws.Range("E20").Value = 1.048313546004334
$ws.Range("F20").Value = 1.053780452855438
$ws.Range("I20").Value = 1.032035200416855
$ws.Range("J20").Value = 1.037490523700523
$ws.Range("K20").Value = 1.036389478658487
$ws.Range("L20").Value = 1.051644624949225
$ws.Range("M20").Value = 1.057093029465751
$ws.Range("N20").Value = 1.038963879246004

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029407222308404
$ws.Range("D21").Value = 1.031649351775939
$ws.Range("E21").Value = 1.046481904640538
$ws.Range("F21").Value = 1.051768569903083
$ws.Range("I21").Value = 1.031583512949368
$ws.Range("J21").Value = 1.036040843415071
$ws.Range("K21").Value = 1.035245756890688
$ws.Range("L21").Value = 1.050023465134133
$ws.Range("M21").Value = 1.055290995477421
$ws.Range("N21").Value = 1.037512140248266

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028237284420058
$ws.Range("D22").Value = 1.030793558252417
$ws.Range("E22").Value = 1.045328243188202
$ws.Range("F22").Value = 1.050501461436864
$ws.Range("I22").Value = 1.031296271821579
$ws.Range("J22").Value = 1.035125610762432
$ws.Range("K22").Value = 1.034522911711704
$ws.Range("L22").Value = 1.049001667134939
$ws.Range("M22").Value = 1.054155392541662
$ws.Range("N22").Value = 1.036595607860281

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028857831779814
$ws.Range("D23").Value = 1.031247533817561
$ws.Range("E23").Value = 1.045940018524227
$ws.Range("F23").Value = 1.051173389227819
$ws.Range("I23").Value = 1.031448852984364
$ws.Range("J23").Value = 1.035611153286835
$ws.Range("K23").Value = 1.034906463332512
$ws.Range("L23").Value = 1.049543583803669
$ws.Range("M23").Value = 1.054757647560098
$ws.Range("N23").Value = 1.037081839910762

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.031293391037474
$ws.Range("D24").Value = 1.033028125402474
$ws.Range("E24").Value = 1.04834422022638
$ws.Range("F24").Value = 1.053814146891567
$ws.Range("I24").Value = 1.032042718280108
$ws.Range("J24").Value = 1.037514764883881
$ws.Range("K24").Value = 1.036408590544932
$ws.Range("L24").Value = 1.051671762296593
$ws.Range("M24").Value = 1.0571231978911
$ws.Range("N24").Value = 1.038988154854623

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034104776961387
$ws.Range("D25").Value = 1.035080941841678
$ws.Range("E25").Value = 1.051125755971845
$ws.Range("F25").Value = 1.05686970968785
$ws.Range("I25").Value = 1.032717865110588
$ws.Range("J25").Value = 1.039707807601275
$ws.Range("K25").Value = 1.03813573263671
$ws.Range("L25").Value = 1.054130887086238
$ws.Range("M25").Value = 1.0598574625193207
$ws.Range("N25").Value = 1.041184311944222
